$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-10, columns B,C,E,F,G,I,J,K,M,O,P,Q
# Columns D,H,L,N are unchanged.

$data = @{
    2  = @{ B=-82.008;             C=529.051; E=0.408; F=-84.502;             G=529.051; I=0.396; J=-27.155; K=529.051; M=1.233; O=113.399; P=51.336; Q=6.614 }
    3  = @{ B=-92.49299999999999;  C=529.051; E=0.362; F=-91.664;             G=529.051; I=0.365; J=-30.627; K=529.051; M=1.093; O=113.399; P=51.336; Q=6.614 }
    4  = @{ B=-92.205;             C=529.051; E=0.363; F=-91.93300000000001;  G=529.051; I=0.364; J=-30.531; K=529.051; M=1.096; O=113.399; P=51.336; Q=6.614 }
    5  = @{ B=-91.66200000000001;  C=529.051; E=0.365; F=-91.681;             G=529.051; I=0.365; J=-30.352; K=529.051; M=1.103; O=113.399; P=51.336; Q=6.614 }
    6  = @{ B=-91.43300000000001;  C=529.051; E=0.366; F=-91.49299999999999;  G=529.051; I=0.366; J=-30.276; K=529.051; M=1.106; O=113.399; P=51.336; Q=6.614 }
    7  = @{ B=-91.66200000000001;  C=529.051; E=0.365; F=-91.681;             G=529.051; I=0.365; J=-30.352; K=529.051; M=1.103; O=113.399; P=51.336; Q=6.614 }
    8  = @{ B=-92.205;             C=529.051; E=0.363; F=-91.93300000000001;  G=529.051; I=0.364; J=-30.531; K=529.051; M=1.096; O=113.399; P=51.336; Q=6.614 }
    9  = @{ B=-92.49299999999999;  C=529.051; E=0.362; F=-91.664;             G=529.051; I=0.365; J=-30.627; K=529.051; M=1.093; O=113.399; P=51.336; Q=6.614 }
    10 = @{ B=-82.008;             C=529.051; E=0.408; F=-84.502;             G=529.051; I=0.396; J=-27.155; K=529.051; M=1.233; O=113.399; P=51.336; Q=6.614 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
